{"js": "// Update the two-digit division answer key table.\n// Each populated row of the table holds 5 \"a\u00f7b=q, r\" answers; we replace\n// the text of specific cells (identified by row/column index) with their\n// new values, leaving every other cell (and all paragraph/run formatting)\n// untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> { column index -> new text }\nconst updates = {\n  0: {\n    0: \"12\u00f73=4, 0\",\n    1: \"49\u00f74=12, 1\",\n    2: \"64\u00f74=16, 0\",\n    3: \"10\u00f73=3, 1\",\n    4: \"23\u00f74=5, 3\",\n  },\n  4: {\n    0: \"15\u00f78=1, 7\",\n    1: \"32\u00f74=8, 0\",\n    2: \"57\u00f79=6, 3\",\n    3: \"96\u00f78=12, 0\",\n    4: \"60\u00f78=7, 4\",\n  },\n  8: {\n    0: \"39\u00f75=7, 4\",\n    1: \"30\u00f77=4, 2\",\n    2: \"83\u00f79=9, 2\",\n    3: \"35\u00f74=8, 3\",\n    4: \"44\u00f73=14, 2\",\n  },\n  12: {\n    0: \"23\u00f74=5, 3\",\n    // column 1 (\"55\u00f72=27, 1\") is unchanged\n    2: \"73\u00f78=9, 1\",\n    3: \"92\u00f79=10, 2\",\n    4: \"64\u00f74=16, 0\",\n  },\n  16: {\n    0: \"15\u00f74=3, 3\",\n    1: \"71\u00f79=7, 8\",\n    2: \"65\u00f75=13, 0\",\n    3: \"63\u00f72=31, 1\",\n    4: \"57\u00f78=7, 1\",\n  },\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const row = Number(rowIndex);\n  const cols = updates[rowIndex];\n  for (const colIndex of Object.keys(cols)) {\n    const col = Number(colIndex);\n    const newText = cols[colIndex];\n\n    const cell = table.getCell(row, col);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    // Replace the text of the cell's (single) paragraph in place so the\n    // existing run/paragraph formatting (fonts, size, justification) is\n    // preserved exactly, instead of replacing the whole cell body.\n    paragraphs.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit division answer key table.\n# The table's populated rows each hold 5 \"a\u00f7b=q, r\" answers; replace the\n# text of specific cells (1-based Row/Column, matching Word's COM Cell()\n# indexer) with new values. Setting Cell.Range.Text only rewrites the\n# run's text, so existing paragraph/run formatting (fonts, size,\n# justification) on that cell is left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @{ Row = 1;  Column = 1; Text = \"12\u00f73=4, 0\" },\n    @{ Row = 1;  Column = 2; Text = \"49\u00f74=12, 1\" },\n    @{ Row = 1;  Column = 3; Text = \"64\u00f74=16, 0\" },\n    @{ Row = 1;  Column = 4; Text = \"10\u00f73=3, 1\" },\n    @{ Row = 1;  Column = 5; Text = \"23\u00f74=5, 3\" },\n\n    @{ Row = 5;  Column = 1; Text = \"15\u00f78=1, 7\" },\n    @{ Row = 5;  Column = 2; Text = \"32\u00f74=8, 0\" },\n    @{ Row = 5;  Column = 3; Text = \"57\u00f79=6, 3\" },\n    @{ Row = 5;  Column = 4; Text = \"96\u00f78=12, 0\" },\n    @{ Row = 5;  Column = 5; Text = \"60\u00f78=7, 4\" },\n\n    @{ Row = 9;  Column = 1; Text = \"39\u00f75=7, 4\" },\n    @{ Row = 9;  Column = 2; Text = \"30\u00f77=4, 2\" },\n    @{ Row = 9;  Column = 3; Text = \"83\u00f79=9, 2\" },\n    @{ Row = 9;  Column = 4; Text = \"35\u00f74=8, 3\" },\n    @{ Row = 9;  Column = 5; Text = \"44\u00f73=14, 2\" },\n\n    @{ Row = 13; Column = 1; Text = \"23\u00f74=5, 3\" },\n    # Row 13, Column 2 (\"55\u00f72=27, 1\") is unchanged.\n    @{ Row = 13; Column = 3; Text = \"73\u00f78=9, 1\" },\n    @{ Row = 13; Column = 4; Text = \"92\u00f79=10, 2\" },\n    @{ Row = 13; Column = 5; Text = \"64\u00f74=16, 0\" },\n\n    @{ Row = 17; Column = 1; Text = \"15\u00f74=3, 3\" },\n    @{ Row = 17; Column = 2; Text = \"71\u00f79=7, 8\" },\n    @{ Row = 17; Column = 3; Text = \"65\u00f75=13, 0\" },\n    @{ Row = 17; Column = 4; Text = \"63\u00f72=31, 1\" },\n    @{ Row = 17; Column = 5; Text = \"57\u00f78=7, 1\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u.Row, $u.Column)\n    $cell.Range.Text = $u.Text\n}\n"}
